# Updated scrum log from yesterday, updated Burndown chart from today
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Yesterday's scrum log update: 2 story points moved from "Planned remaining" (D11)
$ws.Range("D11").Value = 2

# Move the selection to where the user last clicked while reviewing the chart
$ws.Range("D16").Select()
